$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Phase 1: append 5 newly-researched assets under the existing table ---
# (original table is A1:E6, so these land on rows 7-11 first; an insert in
# phase 2 later pushes them down to their final resting rows 9-13)

$ws.Cells.Item(7, 1).Value = "ICSUSPI ID Equity"
$ws.Cells.Item(7, 2).Value = "BlackRock ICS US Dollar Liquidity Fund"
$ws.Cells.Item(7, 6).Value = "Global Fund?"
$ws.Cells.Item(7, 5).Value = "Fixed Income, Cash equivalents"
$ws.Cells.Item(7, 4).Value = "HK"

$ws.Cells.Item(8, 1).Value = "JPMULCD LX Equity"
$ws.Cells.Item(8, 2).Value = "JPMorgan Liquidity Funds"
$ws.Cells.Item(8, 4).Value = "HK"
$ws.Cells.Item(8, 5).Value = "Fixed Income, Cash equivalents"
$ws.Cells.Item(8, 6).Value = "Global Fund?"

$ws.Cells.Item(9, 1).Value = "QQQ US Equity"
$ws.Cells.Item(9, 5).Value = "Fund, Exchange Traded Funds"
$ws.Cells.Item(9, 4).Value = "US"
$ws.Cells.Item(9, 2).Value = "Investco QQQ Trust Series 1"

$ws.Cells.Item(10, 1).Value = "USG8116KAB82"
$ws.Cells.Item(10, 2).Value = "SHNTN 2015-1X B"
$ws.Cells.Item(10, 4).Value = "HK"
$ws.Cells.Item(10, 5).Value = "Fixed income, Asset-Backed, Cash"
$ws.Cells.Item(10, 6).Value = "Country not sure, is it Cash or Synthetic?"

$ws.Cells.Item(11, 1).Value = "SPY US Equity"
$ws.Cells.Item(11, 2).Value = "SPDR S&P 500 ETF Trust"
$ws.Cells.Item(11, 4).Value = "US"
$ws.Cells.Item(11, 5).Value = "Fund, Exchange Traded Funds"

# --- Phase 2: insert 2 more China-ETF rows just under the first asset row ---
$ws.Range("A3:A4").EntireRow.Insert()

$ws.Cells.Item(3, 1).Value = "2823 HK Equity"
$ws.Cells.Item(4, 1).Value = "2828 HK Equity"
$ws.Cells.Item(3, 2).Value = "iShares FTSE A50 China ETF"
$ws.Cells.Item(4, 2).Value = "Hang Seng China Enterprises Index ETF"
$ws.Cells.Item(3, 4).Value = "HK"
$ws.Cells.Item(3, 5).Value = "Fund, Exchange Traded Funds"
$ws.Cells.Item(4, 4).Value = "HK"
$ws.Cells.Item(4, 5).Value = "Fund, Exchange Traded Funds"

# --- Phase 3: flag the rows whose CountryCode/AssetType needs follow-up,
#     highlighting them (and the two trailing helper cells) in yellow ---
$ws.Range("D9").Interior.Color = 65535
$ws.Range("F9").Interior.Color = 65535
$ws.Range("D10").Interior.Color = 65535
$ws.Range("F10").Interior.Color = 65535
$ws.Range("D12").Interior.Color = 65535
$ws.Range("F12:I12").Interior.Color = 65535

# --- Phase 4: cosmetic column-width tweaks made while reviewing the sheet ---
$ws.Columns(1).ColumnWidth = 15.75
$ws.Columns(2).ColumnWidth = 31.67
$ws.Columns(5).ColumnWidth = 28.8

$ws.Range("K11").Select()
